$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price text that can look numeric (e.g. "1.00", "6.30").
# Force them to Text format before assigning so Excel keeps the exact string,
# then reset the style back to Normal so no stray style index is left behind.
function Set-TextCell($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws 'D2' '69.066.19'
$ws.Range('E2').Value = '  -0.03%  '
Set-TextCell $ws 'D3' '3.793.96'
$ws.Range('E3').Value = '  +1.08%  '
$ws.Range('E4').Value = '  -0.01%  '
Set-TextCell $ws 'D5' '600.81'
$ws.Range('E5').Value = '  -0.69%  '
Set-TextCell $ws 'D6' '162.76'
$ws.Range('E6').Value = '  -4.00%  '
Set-TextCell $ws 'D7' '3.792.66'
$ws.Range('E7').Value = '  +1.05%  '
$ws.Range('E8').Value = '  +0.01%  '
Set-TextCell $ws 'D9' '0.536'
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('E10').Value = '  +0.51%  '
Set-TextCell $ws 'D11' '6.30'
$ws.Range('E11').Value = '  -1.13%  '
$ws.Range('E12').Value = '  -1.31%  '
Set-TextCell $ws 'D13' '37.15'
$ws.Range('E13').Value = '  -3.32%  '
$ws.Range('E14').Value = '  -1.84%  '
Set-TextCell $ws 'D15' '4.427.70'
$ws.Range('E15').Value = '  +1.10%  '
Set-TextCell $ws 'D16' '3.816.86'
$ws.Range('E16').Value = '  +1.61%  '
Set-TextCell $ws 'D17' '69.154.09'
$ws.Range('E17').Value = '  +0.12%  '
Set-TextCell $ws 'D18' '7.38'
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell $ws 'D19' '17.32'
$ws.Range('E19').Value = '  +0.94%  '
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell $ws 'D20' '0.114'
$ws.Range('E20').Value = '  -0.43%  '
Set-TextCell $ws 'D21' '11.26'
$ws.Range('E21').Value = '  +3.48%  '
Set-TextCell $ws 'D22' '488.27'
$ws.Range('E22').Value = '  -1.37%  '
$ws.Range('E23').Value = '  -1.43%  '
$ws.Range('E24').Value = '  -2.44%  '
Set-TextCell $ws 'D25' '84.46'
$ws.Range('E25').Value = '  -1.25%  '
$ws.Range('E26').Value = '  -4.35%  '
$ws.Range('E27').Value = '  -1.82%  '
Set-TextCell $ws 'D29' '1.00'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E30').Value = '  -0.89%  '
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('E32').Value = '  -6.19%  '
Set-TextCell $ws 'D33' '3.943.64'
$ws.Range('E33').Value = '  +1.17%  '
Set-TextCell $ws 'D34' '31.80'
$ws.Range('E34').Value = '  -0.62%  '
Set-TextCell $ws 'D35' '3.741.49'
$ws.Range('E35').Value = '  +1.48%  '
$ws.Range('E36').Value = '  -2.37%  '
$ws.Range('B37').Value = 'Mantle'
$ws.Range('C37').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell $ws 'D37' '1.02'
$ws.Range('E37').Value = '  +0.36%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws 'D38' '0.140'
$ws.Range('E38').Value = '  +4.91%  '
Set-TextCell $ws 'D39' '5.89'
$ws.Range('E39').Value = '  -0.38%  '
Set-TextCell $ws 'D40' '0.999'
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('E41').Value = '  -1.15%  '
$ws.Range('E42').Value = '  -1.84%  '
Set-TextCell $ws 'D43' '48.46'
$ws.Range('E43').Value = '  -0.64%  '
$ws.Range('E44').Value = '  -0.41%  '
Set-TextCell $ws 'D45' '418.12'
$ws.Range('E45').Value = '  -5.00%  '
Set-TextCell $ws 'D47' '8.37'
$ws.Range('E47').Value = '  -1.26%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws 'D48' '141.68'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell $ws 'D49' '2.815.17'
$ws.Range('E49').Value = '  +1.16%  '
Set-TextCell $ws 'D50' '1.30'
$ws.Range('E50').Value = '  +4.60%  '
Set-TextCell $ws 'D51' '39.39'
$ws.Range('E51').Value = '  -2.70%  '
